# Auto-generated edit script: updates market-price derived columns (H:N)
# for specific leve rows across multiple crafting-class sheets, matching
# the scheduled market-data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 480.33334
$ws.Range("I6").Value = 77.333336
$ws.Range("J6").Value = 883.3333
$ws.Range("K6").Value = 232.000008
$ws.Range("L6").Value = 2649.9999
$ws.Range("M6").Value = -120.000008
$ws.Range("N6").Value = -2873.9999

$ws.Range("H43").Value = 22892.2
$ws.Range("I43").Value = 111111
$ws.Range("J43").Value = 837.5
$ws.Range("K43").Value = 111111
$ws.Range("L43").Value = 837.5
$ws.Range("M43").Value = -111042
$ws.Range("N43").Value = -975.5

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H94").Value = 1999.5
$ws.Range("I94").Value = 1999.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1999.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1548.5

$ws.Range("H138").Value = 2553.9656
$ws.Range("I138").Value = 1096.2667
$ws.Range("J138").Value = 4115.7856
$ws.Range("K138").Value = 3288.800099999999
$ws.Range("L138").Value = 12347.3568
$ws.Range("M138").Value = 1851.199900000001
$ws.Range("N138").Value = -22627.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1099.75
$ws.Range("I2").Value = 1200
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -1087
$ws.Range("N2").Value = -1225.5

$ws.Range("H4").Value = 594.5
$ws.Range("I4").Value = 190
$ws.Range("J4").Value = 999
$ws.Range("K4").Value = 190
$ws.Range("L4").Value = 999
$ws.Range("M4").Value = -74
$ws.Range("N4").Value = -1231

$ws.Range("H32").Value = 5660.354
$ws.Range("I32").Value = 2251.2246
$ws.Range("J32").Value = 16100.8125
$ws.Range("K32").Value = 2251.2246
$ws.Range("L32").Value = 16100.8125
$ws.Range("M32").Value = -1964.2246
$ws.Range("N32").Value = -16674.8125

$ws.Range("H61").Value = 1761.5526
$ws.Range("I61").Value = 1706.3334
$ws.Range("J61").Value = 1829.7646
$ws.Range("K61").Value = 1706.3334
$ws.Range("L61").Value = 1829.7646
$ws.Range("M61").Value = -1494.3334
$ws.Range("N61").Value = -2253.7646

$ws.Range("H74").Value = 3082.125
$ws.Range("I74").Value = 3457.2666
$ws.Range("J74").Value = 2456.889
$ws.Range("K74").Value = 3457.2666
$ws.Range("L74").Value = 2456.889
$ws.Range("M74").Value = -2583.2666
$ws.Range("N74").Value = -4204.889

$ws.Range("H77").Value = 3082.125
$ws.Range("I77").Value = 3457.2666
$ws.Range("J77").Value = 2456.889
$ws.Range("K77").Value = 17286.333
$ws.Range("L77").Value = 12284.445
$ws.Range("M77").Value = -12918.333
$ws.Range("N77").Value = -21020.445

$ws.Range("H116").Value = 1099.75
$ws.Range("I116").Value = 1200
$ws.Range("J116").Value = 999.5
$ws.Range("K116").Value = 1200
$ws.Range("L116").Value = 999.5
$ws.Range("M116").Value = 1094
$ws.Range("N116").Value = -5587.5

$ws.Range("H132").Value = 1963.2407
$ws.Range("I132").Value = 1559.1538
$ws.Range("J132").Value = 3013.8667
$ws.Range("K132").Value = 4677.4614
$ws.Range("L132").Value = 9041.6001
$ws.Range("M132").Value = -2147.4614
$ws.Range("N132").Value = -14101.6001

$ws.Range("H136").Value = 1761.5526
$ws.Range("I136").Value = 1706.3334
$ws.Range("J136").Value = 1829.7646
$ws.Range("K136").Value = 5119.0002
$ws.Range("L136").Value = 5489.293799999999
$ws.Range("M136").Value = -2569.0002
$ws.Range("N136").Value = -10589.2938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1099.75
$ws.Range("I3").Value = 1200
$ws.Range("J3").Value = 999.5
$ws.Range("K3").Value = 1200
$ws.Range("L3").Value = 999.5
$ws.Range("M3").Value = -1086
$ws.Range("N3").Value = -1227.5

$ws.Range("H94").Value = 11730.895
$ws.Range("I94").Value = 1374.0714
$ws.Range("J94").Value = 40730
$ws.Range("K94").Value = 1374.0714
$ws.Range("L94").Value = 40730
$ws.Range("M94").Value = -923.0714
$ws.Range("N94").Value = -41632

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 34
$ws.Range("I7").Value = 34
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 34
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 79
$ws.Range("N7").ClearContents()

$ws.Range("H105").Value = 959
$ws.Range("I105").Value = 734
$ws.Range("J105").Value = 1334
$ws.Range("K105").Value = 734
$ws.Range("L105").Value = 1334
$ws.Range("M105").Value = 1013

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 260.4
$ws.Range("I11").Value = 166.66667
$ws.Range("J11").Value = 401
$ws.Range("K11").Value = 500.00001
$ws.Range("L11").Value = 1203
$ws.Range("M11").Value = -360.00001
$ws.Range("N11").Value = -1483

$ws.Range("H62").Value = 6118
$ws.Range("I62").Value = 500
$ws.Range("J62").Value = 6679.8
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 20039.4
$ws.Range("M62").Value = -814
$ws.Range("N62").Value = -21411.4

$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("K63").Value = 90000
$ws.Range("L63").Value = 90000
$ws.Range("M63").Value = -89251
$ws.Range("N63").Value = -91498

$ws.Range("H64").Value = 4278.826
$ws.Range("I64").Value = 1979.8
$ws.Range("J64").Value = 4917.4443
$ws.Range("K64").Value = 5939.4
$ws.Range("L64").Value = 14752.3329
$ws.Range("M64").Value = -5669.4
$ws.Range("N64").Value = -15292.3329

$ws.Range("H65").Value = 6118
$ws.Range("I65").Value = 500
$ws.Range("J65").Value = 6679.8
$ws.Range("K65").Value = 4500
$ws.Range("L65").Value = 60118.2
$ws.Range("M65").Value = -1068
$ws.Range("N65").Value = -66982.20000000001

$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("K66").Value = 270000
$ws.Range("L66").Value = 270000
$ws.Range("M66").Value = -266256
$ws.Range("N66").Value = -277488

$ws.Range("H67").Value = 4278.826
$ws.Range("I67").Value = 1979.8
$ws.Range("J67").Value = 4917.4443
$ws.Range("K67").Value = 5939.4
$ws.Range("L67").Value = 14752.3329
$ws.Range("M67").Value = -5003.4
$ws.Range("N67").Value = -16624.3329

$ws.Range("H69").Value = 2506
$ws.Range("I69").Value = 2506
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 7518
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -6707
$ws.Range("N69").ClearContents()

$ws.Range("H70").Value = 5966.6665
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 6372.727
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 19118.181
$ws.Range("M70").Value = -4185
$ws.Range("N70").Value = -19748.181

$ws.Range("H72").Value = 2506
$ws.Range("I72").Value = 2506
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 22554
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -18498
$ws.Range("N72").ClearContents()

$ws.Range("H73").Value = 5966.6665
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 6372.727
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 19118.181
$ws.Range("M73").Value = -3408
$ws.Range("N73").Value = -21302.181

$ws.Range("H75").Value = 3000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 9000
$ws.Range("N75").Value = -10996
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 3000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 27000
$ws.Range("N78").Value = -36984
$ws.Range("M78").ClearContents()

$ws.Range("H131").Value = 899.1279
$ws.Range("I131").Value = 587.65
$ws.Range("J131").Value = 993.51514
$ws.Range("K131").Value = 1762.95
$ws.Range("L131").Value = 2980.54542
$ws.Range("M131").Value = 3277.05
$ws.Range("N131").Value = -13060.54542

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1625.8536
$ws.Range("I132").Value = 1066.9667
$ws.Range("J132").Value = 3150.0908
$ws.Range("K132").Value = 3200.9001
$ws.Range("L132").Value = 9450.2724
$ws.Range("M132").Value = -670.9000999999998
$ws.Range("N132").Value = -14510.2724

$ws.Range("H138").Value = 26322.223
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 26322.223
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 26322.223
$ws.Range("N138").Value = -36602.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3893.88
$ws.Range("I40").Value = 4020.5293
$ws.Range("J40").Value = 3624.75
$ws.Range("K40").Value = 4020.5293
$ws.Range("L40").Value = 3624.75
$ws.Range("M40").Value = -3884.5293

$ws.Range("H122").Value = 3097.9167
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3097.9167
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9293.750100000001
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 786.8570999999999
$ws.Range("I96").Value = 655
$ws.Range("J96").Value = 962.6667
$ws.Range("K96").Value = 655
$ws.Range("L96").Value = 962.6667
$ws.Range("M96").Value = 718
$ws.Range("N96").Value = -3708.6667

$ws.Range("H100").Value = 6757.7646
$ws.Range("I100").Value = 12925.25
$ws.Range("J100").Value = 1275.5555
$ws.Range("K100").Value = 25850.5
$ws.Range("L100").Value = 2551.111
$ws.Range("M100").Value = -25309.5
$ws.Range("N100").Value = -3633.111
